$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15's L-column comment gets a run-count prefix added to its existing text.
$ws.Range("L15").Value = "8 runs for 5%. New normalisation meathod. Ideal (coherent, noiseless) Interferogram peaks are now set to 1."

# New "Comments" entries added to the L column (detector table) for rows 13-21,
# giving each detector row a width/run-count note.
$ws.Range("L13").Value = "484 runs for 5%"
$ws.Range("L14").Value = "27 runs for 5%"
$ws.Range("L16").Value = "15 runs for 5%"
$ws.Range("L17").Value = "324 runs for 5%"
$ws.Range("L18").Value = "501 runs for 5%"
$ws.Range("L19").Value = "501 runs for 5%"
$ws.Range("L20").Value = "218 runs for 5%"
$ws.Range("L21").Value = "501 runs for 5%"

$ws.Range("L22").Select()
